$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 previously held "Einstellung Kamera" in A24 (an orphan leftover
# row, disconnected from the A20:B23 "Beinflussende Parameter" block).
# Remove that stray entry entirely (clearing the cell also drops the now
# -unused shared string, shifting later shared-string indices down by one
# automatically).
$ws.Range("A24").ClearContents() | Out-Null

# Give the now content-less B24 a plain "finished/boxed" placeholder look:
# solid white (background 1) fill plus a thin box border all the way
# round, matching the rest of the sheet's manual formatting conventions.
$b24 = $ws.Range("B24")
$b24.Interior.Color = 255
$b24.Interior.ThemeColor = 2
$b24.Borders.LineStyle = 1
$b24.Borders.Weight = 2

# Move the active selection from the old B27 spot to A24.
$ws.Range("A24").Select() | Out-Null
